$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 16
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = -3
$ws.Range("F9").Value = -11
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 2
